# "small adds to Chapter 3":
#  - add a new "Sheet2" after the existing "Sheet1" and make it the active sheet
#  - populate it with an OUTLIERS table (header + 17 data rows)
#  - style the table: centered text, and a centered short-date for the Date column
#  - autofit column A, and leave A2:C19 selected (matches the saved sheet view)

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 (keeps Sheet1 first / Sheet2 second,
# and Excel makes the newly inserted sheet the active one).
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)

# Title row
$ws.Range("A1").Value = "OUTLIERS"

# Header row
$ws.Range("A2").Value = "Date"
$ws.Range("B2").Value = "Sample"
$ws.Range("C2").Value = "N form"

$dates = @(42296, 41508, 42115, 41753, 41844, 41844, 41810, 42058, 41483, 41508, 42058, 41025, 42172, 42206, 41550, 42543, 42576)
$samples = @("T11T2", "T11T2", "T12T1", "T12T1", "T10S2", "T10S2", "T12S2", "T10S1", "T10S2", "C10S1", "T10S2", "T10S2", "C31D1", "C31D1", "C31D1", "C31D1", "C30D1")
$nforms = @("NO3.N", "NH4.N", "NH4.N", "NH4.N", "NO3.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N", "NH4.N")

# Write column by column (matches how the shared-string table ends up ordered)
for ($i = 0; $i -lt $dates.Length; $i = $i + 1) {
    $ws.Cells.Item(3 + $i, 1).Value = $dates[$i]
}
for ($i = 0; $i -lt $samples.Length; $i = $i + 1) {
    $ws.Cells.Item(3 + $i, 2).Value = $samples[$i]
}
for ($i = 0; $i -lt $nforms.Length; $i = $i + 1) {
    $ws.Cells.Item(3 + $i, 3).Value = $nforms[$i]
}

# Apply the centered style to the whole table (A2:C19) from a single styled
# scratch cell, so every cell in the block shares ONE cellXf entry.
$ws.Range("E1").Value = 0
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").Copy($ws.Range("A2:C19")) | Out-Null

# Apply the centered short-date style to the Date column (A3:A19), again
# from a single styled scratch cell so it ends up as one shared cellXf.
$ws.Range("E2").Value = 1
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Copy($ws.Range("A3:A19")) | Out-Null

# Drop the scratch cells
$ws.Range("E1:E2").Clear() | Out-Null

# The Copy() broadcasts above clobbered the cell values with the scratch
# values, so re-enter them now that the formatting is in place.
$ws.Range("A2").Value = "Date"
$ws.Range("B2").Value = "Sample"
$ws.Range("C2").Value = "N form"
for ($i = 0; $i -lt $dates.Length; $i = $i + 1) {
    $ws.Cells.Item(3 + $i, 1).Value = $dates[$i]
}
for ($i = 0; $i -lt $samples.Length; $i = $i + 1) {
    $ws.Cells.Item(3 + $i, 2).Value = $samples[$i]
}
for ($i = 0; $i -lt $nforms.Length; $i = $i + 1) {
    $ws.Cells.Item(3 + $i, 3).Value = $nforms[$i]
}

# Size column A to fit its (now date-formatted) contents
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the table selected, matching the saved sheet view
$ws.Range("A2:C19").Select() | Out-Null
